# iNatStruct.xlsx edit: recompute dominant-color percentage/hex ordering
# and refresh the derived RGB triplet values (R1/G1/B1, R2/G2/B2, R3/G3/B3)
# after removing the now-unused python color-clustering library.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 126.5631645844316
$ws.Range("K2").Value = 135.7539758812044
$ws.Range("L2").Value = 110.2767151315988
$ws.Range("M2").Value = 61.49758395955317
$ws.Range("N2").Value = 64.65939848329957
$ws.Range("O2").Value = 39.38283493542211
$ws.Range("P2").Value = 183.6292628753858
$ws.Range("Q2").Value = 189.6336279363796
$ws.Range("R2").Value = 180.0335631545592
$ws.Range("S2").Value = "[0.34, 0.31, 0.35]"
$ws.Range("T2").Value = "['#7f886e', '#3d4127', '#b8beb4']"

# Row 3
$ws.Range("J3").Value = 116.8762323920106
$ws.Range("K3").Value = 130.04179096491
$ws.Range("L3").Value = 80.58489960269071
$ws.Range("M3").Value = 53.22184001949568
$ws.Range("N3").Value = 64.81027934241865
$ws.Range("O3").Value = 28.92279027103484
$ws.Range("P3").Value = 191.1821714445269
$ws.Range("Q3").Value = 200.2241948994773
$ws.Range("R3").Value = 137.4939966492951
$ws.Range("S3").Value = "[0.4, 0.37, 0.23]"
$ws.Range("T3").Value = "['#758251', '#35411d', '#bfc889']"

# Row 4
$ws.Range("J4").Value = 196.1845401388855
$ws.Range("K4").Value = 187.4492326503656
$ws.Range("L4").Value = 155.6496841247063
$ws.Range("M4").Value = 61.90421877825159
$ws.Range("N4").Value = 59.49599815932039
$ws.Range("O4").Value = 44.64111624238532
$ws.Range("P4").Value = 130.0051064235949
$ws.Range("Q4").Value = 119.6158251407843
$ws.Range("R4").Value = 97.71941713594822
$ws.Range("S4").Value = "[0.34, 0.32, 0.34]"
$ws.Range("T4").Value = "['#c4bb9c', '#3e3b2d', '#827862']"

# Row 5
$ws.Range("J5").Value = 191.3518264311085
$ws.Range("K5").Value = 206.0431961207371
$ws.Range("L5").Value = 221.9472304922189
$ws.Range("M5").Value = 48.64646084337149
$ws.Range("N5").Value = 76.94636237256844
$ws.Range("O5").Value = 17.71605653384222
$ws.Range("P5").Value = 97.96231550041296
$ws.Range("Q5").Value = 128.9413769430669
$ws.Range("R5").Value = 56.08100818909594
$ws.Range("S5").Value = "[0.6, 0.18, 0.22]"
$ws.Range("T5").Value = "['#bfcede', '#314d12', '#628138']"

# Row 6
$ws.Range("J6").Value = 187.5822362617607
$ws.Range("K6").Value = 202.1249023347274
$ws.Range("L6").Value = 220.0308160196237
$ws.Range("M6").Value = 41.70808148114206
$ws.Range("N6").Value = 72.91859026225309
$ws.Range("O6").Value = 19.47520346483464
$ws.Range("P6").Value = 162.1814435328357
$ws.Range("Q6").Value = 176.7655866006281
$ws.Range("R6").Value = 187.6934477674003

# Row 7
$ws.Range("J7").Value = 204.4227235919868
$ws.Range("K7").Value = 200.2727221958501
$ws.Range("L7").Value = 192.8627732946076
$ws.Range("M7").Value = 39.58090791179841
$ws.Range("N7").Value = 54.97963683527826
$ws.Range("O7").Value = 20.33447470817669
$ws.Range("P7").Value = 101.9141017019259
$ws.Range("Q7").Value = 119.8082225440996
$ws.Range("R7").Value = 77.92479695826066
$ws.Range("S7").Value = "[0.38, 0.21, 0.41]"
$ws.Range("T7").Value = "['#ccc8c1', '#283714', '#66784e']"

# Row 8
$ws.Range("J8").Value = 183.9214150159144
$ws.Range("K8").Value = 199.0223889960659
$ws.Range("L8").Value = 219.7645508869342
$ws.Range("M8").Value = 92.38540466589484
$ws.Range("N8").Value = 107.7785858294953
$ws.Range("O8").Value = 68.48819124424988
$ws.Range("P8").Value = 54.59066633546077
$ws.Range("Q8").Value = 49.60212382786062
$ws.Range("R8").Value = 22.57256411849005
$ws.Range("S8").Value = "[0.68, 0.15, 0.17]"
$ws.Range("T8").Value = "['#b8c7dc', '#5c6c44', '#373217']"

# Row 9
$ws.Range("J9").Value = 24.62294724073291
$ws.Range("K9").Value = 39.45588841721808
$ws.Range("L9").Value = 24.08240145545861
$ws.Range("M9").Value = 172.6558463915306
$ws.Range("N9").Value = 187.4755683072177
$ws.Range("O9").Value = 210.0098256455701
$ws.Range("P9").Value = 153.0578181818192
$ws.Range("Q9").Value = 167.0655757575737
$ws.Range("R9").Value = 187.5536969696962

# Row 10
$ws.Range("J10").Value = 171.5424700663058
$ws.Range("K10").Value = 185.9211849449001
$ws.Range("L10").Value = 209.7348184897223
$ws.Range("M10").Value = 62.50506083223027
$ws.Range("N10").Value = 61.79112565176095
$ws.Range("O10").Value = 41.62810551067622
$ws.Range("P10").Value = 150.8932243397807
$ws.Range("Q10").Value = 163.7523723523916
$ws.Range("R10").Value = 185.3814140646039
$ws.Range("S10").Value = "[0.56, 0.1, 0.34]"
$ws.Range("T10").Value = "['#acbad2', '#3f3e2a', '#97a4b9']"

# Row 11
$ws.Range("J11").Value = 184.8598147768016
$ws.Range("K11").Value = 184.1953174996002
$ws.Range("L11").Value = 181.3286926912922
$ws.Range("M11").Value = 141.6376205941343
$ws.Range("N11").Value = 141.5192234215295
$ws.Range("O11").Value = 136.4505802846509
$ws.Range("P11").Value = 112.4421615939599
$ws.Range("Q11").Value = 76.14897988436529
$ws.Range("R11").Value = 62.99483969612039
$ws.Range("S11").Value = "[0.44, 0.45, 0.11]"
$ws.Range("T11").Value = "['#b9b8b5', '#8e8e88', '#704c3f']"
